# NYPD CompStat 70th Precinct weekly report — refreshed with new crime data
# (volume/week-range header text + all Week-to-Date / 28-Day / Year-to-Date /
# 2-Year / 13-Year / 30-Year figures on the crime-category rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: "Volume 30   Number  35" -> "...Number  36"
# and "Report Covering the Week  8/28/2023  Through  9/3/2023"
#   -> "...9/4/2023  Through  9/10/2023"
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/4/2023  Through  9/10/2023"

# ---------------------------------------------------------------------------
# Helper: turn a numeric cell into the literal text "0" (used where a
# category had data previously but now has none for the Week-to-Date column)
# while preserving the existing "General" / text style already used by
# other zero-count cells on the sheet (style copied from C14).
# ---------------------------------------------------------------------------
function Set-TextZero($addr) {
    $ws.Range($addr).Value = "'0"
    $ws.Range("C14").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------------
$ws.Range("N14").Value = -87.878787878787

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
Set-TextZero "C15"
$ws.Range("E15").Value = -100
$ws.Range("J15").Value = 24
$ws.Range("K15").Value = -41.666666666666
$ws.Range("M15").Value = -22.222222222222
$ws.Range("N15").Value = -84.782608695652

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 7
$ws.Range("H16").Value = -53.333333333333
$ws.Range("I16").Value = 105
$ws.Range("J16").Value = 142
$ws.Range("K16").Value = -26.056338028169
$ws.Range("L16").Value = -6.25
$ws.Range("M16").Value = -59.770114942528
$ws.Range("N16").Value = -92.934051144010

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -37.5
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = -34.482758620689
$ws.Range("I17").Value = 227
$ws.Range("J17").Value = 259
$ws.Range("K17").Value = -12.355212355212
$ws.Range("L17").Value = 2.714932126696
$ws.Range("M17").Value = -13.358778625954
$ws.Range("N17").Value = -63.738019169329

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 84
$ws.Range("J18").Value = 122
$ws.Range("K18").Value = -31.147540983606
$ws.Range("L18").Value = -33.333333333333
$ws.Range("M18").Value = -59.420289855072
$ws.Range("N18").Value = -95.763993948562

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -42.857142857142
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 50
$ws.Range("H19").Value = -36
$ws.Range("I19").Value = 357
$ws.Range("J19").Value = 372
$ws.Range("K19").Value = -4.032258064516
$ws.Range("L19").Value = 9.509202453987
$ws.Range("M19").Value = -19.230769230769
$ws.Range("N19").Value = -52.652519893899

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 12
$ws.Range("H20").Value = 71.428571428571
$ws.Range("I20").Value = 77
$ws.Range("J20").Value = 86
$ws.Range("K20").Value = -10.465116279069
$ws.Range("L20").Value = -17.204301075268
$ws.Range("M20").Value = -43.795620437956
$ws.Range("N20").Value = -95.546558704453

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -36.666666666666
$ws.Range("F21").Value = 82
$ws.Range("G21").Value = 112
$ws.Range("H21").Value = -26.785714285714
$ws.Range("I21").Value = 868
$ws.Range("J21").Value = 1006
$ws.Range("K21").Value = -13.717693836978
$ws.Range("L21").Value = -3.448275862068
$ws.Range("M21").Value = -34.932533733133
$ws.Range("N21").Value = -87.050574369685

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -75
$ws.Range("J22").Value = 10
$ws.Range("K22").Value = 60
$ws.Range("L22").Value = 77.777777777777
$ws.Range("M22").Value = -23.809523809523

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 35
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 111
$ws.Range("G24").Value = 155
$ws.Range("H24").Value = -28.387096774193
$ws.Range("I24").Value = 1255
$ws.Range("J24").Value = 974
$ws.Range("K24").Value = 28.850102669404
$ws.Range("L24").Value = 40.695067264574
$ws.Range("M24").Value = 40.380313199105

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 44
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = 41.935483870967
$ws.Range("I25").Value = 422
$ws.Range("J25").Value = 399
$ws.Range("K25").Value = 5.764411027568
$ws.Range("L25").Value = 16.253443526170
$ws.Range("M25").Value = -21.851851851851

# ---------------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------------
Set-TextZero "C26"
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 4
$ws.Range("I26").Value = 24
$ws.Range("J26").Value = 35
$ws.Range("K26").Value = -31.428571428571
$ws.Range("L26").Value = -27.272727272727

# ---------------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------------
Set-TextZero "C27"
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -25
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = -10
$ws.Range("L27").Value = -2.173913043478

# ---------------------------------------------------------------------------
# Row 28 - Shooting Vic.
# ---------------------------------------------------------------------------
$ws.Range("N28").Value = -91.428571428571

# ---------------------------------------------------------------------------
# Row 29 - Shooting Inc.
# ---------------------------------------------------------------------------
$ws.Range("N29").Value = -90.588235294117

# ---------------------------------------------------------------------------
# Row 30 - Hate Crimes
# ---------------------------------------------------------------------------
Set-TextZero "F30"
$ws.Range("H30").Value = -100
